$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-30 Tuesday" "2025-10-01 Wednesday"

Replace-Text "473÷4=" "828÷3="
Replace-Text "234÷3=" "269÷9="
Replace-Text "588÷7=" "699÷7="
Replace-Text "942÷5=" "195÷4="
Replace-Text "134÷9=" "263÷2="

Replace-Text "153÷3=" "992÷3="
Replace-Text "478÷4=" "129÷6="
Replace-Text "909÷6=" "325÷8="
Replace-Text "875÷9=" "250÷7="
Replace-Text "627÷5=" "272÷4="

Replace-Text "599÷3=" "436÷4="
Replace-Text "495÷8=" "526÷9="
Replace-Text "232÷7=" "527÷6="
Replace-Text "804÷5=" "598÷7="
Replace-Text "723÷9=" "683÷4="

Replace-Text "907÷6=" "739÷7="
Replace-Text "909÷9=" "654÷7="
Replace-Text "104÷5=" "894÷3="
Replace-Text "927÷6=" "735÷9="
Replace-Text "193÷9=" "107÷9="

Replace-Text "799÷6=" "646÷4="
Replace-Text "971÷6=" "498÷5="
Replace-Text "668÷4=" "221÷8="
Replace-Text "600÷8=" "930÷2="
Replace-Text "817÷5=" "881÷3="
